$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-17 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-18 Thursday", 2) | Out-Null
$d.Content.Find.Execute("262×5=1310", $true, $false, $false, $false, $false, $true, 1, $false, "805×8=6440", 2) | Out-Null
$d.Content.Find.Execute("226×3=678", $true, $false, $false, $false, $false, $true, 1, $false, "685×6=4110", 2) | Out-Null
$d.Content.Find.Execute("111×6=666", $true, $false, $false, $false, $false, $true, 1, $false, "656×3=1968", 2) | Out-Null
$d.Content.Find.Execute("187×6=1122", $true, $false, $false, $false, $false, $true, 1, $false, "784×3=2352", 2) | Out-Null
$d.Content.Find.Execute("903×6=5418", $true, $false, $false, $false, $false, $true, 1, $false, "296×8=2368", 2) | Out-Null
$d.Content.Find.Execute("300×3=900", $true, $false, $false, $false, $false, $true, 1, $false, "434×8=3472", 2) | Out-Null
$d.Content.Find.Execute("789×6=4734", $true, $false, $false, $false, $false, $true, 1, $false, "630×9=5670", 2) | Out-Null
$d.Content.Find.Execute("896×2=1792", $true, $false, $false, $false, $false, $true, 1, $false, "195×9=1755", 2) | Out-Null
$d.Content.Find.Execute("990×5=4950", $true, $false, $false, $false, $false, $true, 1, $false, "252×7=1764", 2) | Out-Null
$d.Content.Find.Execute("934×8=7472", $true, $false, $false, $false, $false, $true, 1, $false, "706×7=4942", 2) | Out-Null
$d.Content.Find.Execute("676×3=2028", $true, $false, $false, $false, $false, $true, 1, $false, "153×5=765", 2) | Out-Null
$d.Content.Find.Execute("389×3=1167", $true, $false, $false, $false, $false, $true, 1, $false, "701×7=4907", 2) | Out-Null
$d.Content.Find.Execute("992×7=6944", $true, $false, $false, $false, $false, $true, 1, $false, "832×8=6656", 2) | Out-Null
$d.Content.Find.Execute("584×3=1752", $true, $false, $false, $false, $false, $true, 1, $false, "763×7=5341", 2) | Out-Null
$d.Content.Find.Execute("194×5=970", $true, $false, $false, $false, $false, $true, 1, $false, "632×3=1896", 2) | Out-Null
$d.Content.Find.Execute("647×9=5823", $true, $false, $false, $false, $false, $true, 1, $false, "770×8=6160", 2) | Out-Null
$d.Content.Find.Execute("955×2=1910", $true, $false, $false, $false, $false, $true, 1, $false, "412×2=824", 2) | Out-Null
$d.Content.Find.Execute("489×8=3912", $true, $false, $false, $false, $false, $true, 1, $false, "657×7=4599", 2) | Out-Null
$d.Content.Find.Execute("482×5=2410", $true, $false, $false, $false, $false, $true, 1, $false, "132×5=660", 2) | Out-Null
$d.Content.Find.Execute("528×4=2112", $true, $false, $false, $false, $false, $true, 1, $false, "945×7=6615", 2) | Out-Null
$d.Content.Find.Execute("483×2=966", $true, $false, $false, $false, $false, $true, 1, $false, "770×9=6930", 2) | Out-Null
$d.Content.Find.Execute("258×4=1032", $true, $false, $false, $false, $false, $true, 1, $false, "943×5=4715", 2) | Out-Null
$d.Content.Find.Execute("600×4=2400", $true, $false, $false, $false, $false, $true, 1, $false, "937×7=6559", 2) | Out-Null
$d.Content.Find.Execute("822×6=4932", $true, $false, $false, $false, $false, $true, 1, $false, "513×4=2052", 2) | Out-Null
$d.Content.Find.Execute("128×2=256", $true, $false, $false, $false, $false, $true, 1, $false, "980×8=7840", 2) | Out-Null
